# Deg to Turn Measurements.xlsx - edit script
#
# Changes applied (per the commit's OOXML diff):
#   - I2: Pulse On Time (us) value 235 -> 190
#     (L2 = K2*I2 recalculates automatically from 5222170 -> 4222180)
#   - Active selection on Sheet1 moves from E15 to L2
#
# (Workbook/window chrome such as fileVersion, calcPr@calcId, the
# bookViews pixel geometry, and the worksheet's pageSetup DPI attributes
# are last-saved-application stamps that Excel itself rewrites on every
# save and are not meaningful, settable document content - they are left
# alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Pulse On Time input; the dependent L2 formula (=K2*I2)
# recalculates automatically.
$ws.Range("I2").Value = 190

# Make sure the sheet is active, then move the selection to L2, matching
# the saved <selection activeCell="L2" sqref="L2"/>.
$ws.Activate()
$ws.Range("L2").Select()
